# Apply the edits described by the commit:
# "created function for adding index, added exceptions for exiting when
#  entered arguments causing errors. Created file and directory for tests."
#
# Concretely this fills in rows 2-5 (columns A-D) of the sheet with
# sample data (a name, an integer, a decimal and a date), moves the
# active selection from E8 to F8, and sets up the page for printing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = "dsanjk"
$ws.Range("B2").Value = 123
$ws.Range("C2").Value = 123.4
$ws.Range("D2").Value = "12/12/2012"

# --- Row 3 ---
$ws.Range("A3").Value = "dsjai"
$ws.Range("B3").Value = 312432
$ws.Range("C3").Value = 412.6
$ws.Range("D3").Value = "04/12/1998"

# --- Row 4 ---
$ws.Range("A4").Value = "asdjkl"
$ws.Range("B4").Value = 532
$ws.Range("C4").Value = 412.42099999999999
$ws.Range("D4").Value = "05/18/2024"

# --- Row 5 ---
$ws.Range("A5").Value = "rjio"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 80.099999999999994
$ws.Range("D5").Value = "01/01/2001"

# Move the active selection from E8 to F8
[void]$ws.Range("F8").Select()

# Configure the page for printing (paper size 9 = A4, portrait orientation)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
